$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Octubre de 2020 a las 13:53"

# Row 16: Iran
$ws.Cells.Item(16,1).Value = "Iran"
$ws.Cells.Item(16,2).Value = 545286
$ws.Cells.Item(16,3).Value = 5616
$ws.Cells.Item(16,4).Value = 438709
$ws.Cells.Item(16,5).Value = 75231
$ws.Cells.Item(16,6).Value = 0
$ws.Cells.Item(16,7).Value = 312
$ws.Cells.Item(16,8).Value = 31346

# Row 20: Banglades
$ws.Cells.Item(20,1).Value = "Banglades"
$ws.Cells.Item(20,2).Value = 393131
$ws.Cells.Item(20,3).Value = 1545
$ws.Cells.Item(20,4).Value = 308845
$ws.Cells.Item(20,5).Value = 78563
$ws.Cells.Item(20,6).Value = 0
$ws.Cells.Item(20,7).Value = 24
$ws.Cells.Item(20,8).Value = 5723

# Row 37: Nepal
$ws.Cells.Item(37,1).Value = "Nepal"
$ws.Cells.Item(37,2).Value = 144872
$ws.Cells.Item(37,3).Value = 5743
$ws.Cells.Item(37,4).Value = 99605
$ws.Cells.Item(37,5).Value = 44476
$ws.Cells.Item(37,6).Value = 0
$ws.Cells.Item(37,7).Value = 26
$ws.Cells.Item(37,8).Value = 791

# Row 38: Bolivia
$ws.Cells.Item(38,1).Value = "Bolivia"
$ws.Cells.Item(38,2).Value = 140037
$ws.Cells.Item(38,3).Value = 147
$ws.Cells.Item(38,4).Value = 105526
$ws.Cells.Item(38,5).Value = 25985
$ws.Cells.Item(38,6).Value = 0
$ws.Cells.Item(38,7).Value = 24
$ws.Cells.Item(38,8).Value = 8526

# Row 42: Emiratos Arabes Unidos
$ws.Cells.Item(42,1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(42,2).Value = 119132
$ws.Cells.Item(42,3).Value = 1538
$ws.Cells.Item(42,4).Value = 111814
$ws.Cells.Item(42,5).Value = 6846
$ws.Cells.Item(42,6).Value = 0
$ws.Cells.Item(42,7).Value = 2
$ws.Cells.Item(42,8).Value = 472

# Row 43: Kuwait
$ws.Cells.Item(43,1).Value = "Kuwait"
$ws.Cells.Item(43,2).Value = 117718
$ws.Cells.Item(43,3).Value = 0
$ws.Cells.Item(43,4).Value = 109198
$ws.Cells.Item(43,5).Value = 7806
$ws.Cells.Item(43,6).Value = 0
$ws.Cells.Item(43,7).Value = 0
$ws.Cells.Item(43,8).Value = 714

# Row 44: Oman
$ws.Cells.Item(44,1).Value = "Oman"
$ws.Cells.Item(44,2).Value = 111484
$ws.Cells.Item(44,3).Value = 451
$ws.Cells.Item(44,4).Value = 97367
$ws.Cells.Item(44,5).Value = 12980
$ws.Cells.Item(44,6).Value = 0
$ws.Cells.Item(44,7).Value = 15
$ws.Cells.Item(44,8).Value = 1137

# Row 82: Bosnia y Herzegovina
$ws.Cells.Item(82,1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(82,2).Value = 36315
$ws.Cells.Item(82,3).Value = 926
$ws.Cells.Item(82,4).Value = 25779
$ws.Cells.Item(82,5).Value = 9505
$ws.Cells.Item(82,6).Value = 0
$ws.Cells.Item(82,7).Value = 14
$ws.Cells.Item(82,8).Value = 1031

# Row 91: Malasia
$ws.Cells.Item(91,1).Value = "Malasia"
$ws.Cells.Item(91,2).Value = 22957
$ws.Cells.Item(91,3).Value = 732
$ws.Cells.Item(91,4).Value = 14931
$ws.Cells.Item(91,5).Value = 7827
$ws.Cells.Item(91,6).Value = 0
$ws.Cells.Item(91,7).Value = 6
$ws.Cells.Item(91,8).Value = 199

# Row 97: Noruega
$ws.Cells.Item(97,1).Value = "Noruega"
$ws.Cells.Item(97,2).Value = 16772
$ws.Cells.Item(97,3).Value = 0
$ws.Cells.Item(97,4).Value = 11863
$ws.Cells.Item(97,5).Value = 4630
$ws.Cells.Item(97,6).Value = 0
$ws.Cells.Item(97,7).Value = 1
$ws.Cells.Item(97,8).Value = 279

# Row 98: Eslovenia
$ws.Cells.Item(98,1).Value = "Eslovenia"
$ws.Cells.Item(98,2).Value = 15982
$ws.Cells.Item(98,3).Value = 1503
$ws.Cells.Item(98,4).Value = 6922
$ws.Cells.Item(98,5).Value = 8860
$ws.Cells.Item(98,6).Value = 0
$ws.Cells.Item(98,7).Value = 8
$ws.Cells.Item(98,8).Value = 200

# Row 99: Zambia
$ws.Cells.Item(99,1).Value = "Zambia"
$ws.Cells.Item(99,2).Value = 15982
$ws.Cells.Item(99,3).Value = 0
$ws.Cells.Item(99,4).Value = 15038
$ws.Cells.Item(99,5).Value = 598
$ws.Cells.Item(99,6).Value = 0
$ws.Cells.Item(99,7).Value = 0
$ws.Cells.Item(99,8).Value = 346

# Row 100: Montenegro
$ws.Cells.Item(100,1).Value = "Montenegro"
$ws.Cells.Item(100,2).Value = 15892
$ws.Cells.Item(100,3).Value = 0
$ws.Cells.Item(100,4).Value = 11581
$ws.Cells.Item(100,5).Value = 4064
$ws.Cells.Item(100,6).Value = 0
$ws.Cells.Item(100,7).Value = 0
$ws.Cells.Item(100,8).Value = 247

# Row 101: Senegal
$ws.Cells.Item(101,1).Value = "Senegal"
$ws.Cells.Item(101,2).Value = 15484
$ws.Cells.Item(101,3).Value = 25
$ws.Cells.Item(101,4).Value = 13975
$ws.Cells.Item(101,5).Value = 1188
$ws.Cells.Item(101,6).Value = 0
$ws.Cells.Item(101,7).Value = 1
$ws.Cells.Item(101,8).Value = 321

# Row 102: Finlandia
$ws.Cells.Item(102,1).Value = "Finlandia"
$ws.Cells.Item(102,2).Value = 14071
$ws.Cells.Item(102,3).Value = 222
$ws.Cells.Item(102,4).Value = 9800
$ws.Cells.Item(102,5).Value = 3916
$ws.Cells.Item(102,6).Value = 0
$ws.Cells.Item(102,7).Value = 4
$ws.Cells.Item(102,8).Value = 355

# Row 136: Malta
$ws.Cells.Item(136,1).Value = "Malta"
$ws.Cells.Item(136,2).Value = 5026
$ws.Cells.Item(136,3).Value = 155
$ws.Cells.Item(136,4).Value = 3331
$ws.Cells.Item(136,5).Value = 1649
$ws.Cells.Item(136,6).Value = 0
$ws.Cells.Item(136,7).Value = 0
$ws.Cells.Item(136,8).Value = 46

# Row 137: Ruanda
$ws.Cells.Item(137,1).Value = "Ruanda"
$ws.Cells.Item(137,2).Value = 4996
$ws.Cells.Item(137,3).Value = 0
$ws.Cells.Item(137,4).Value = 4797
$ws.Cells.Item(137,5).Value = 165
$ws.Cells.Item(137,6).Value = 0
$ws.Cells.Item(137,7).Value = 0
$ws.Cells.Item(137,8).Value = 34

# Row 138: Reunion
$ws.Cells.Item(138,1).Value = "Reunion"
$ws.Cells.Item(138,2).Value = 4921
$ws.Cells.Item(138,3).Value = 0
$ws.Cells.Item(138,4).Value = 4445
$ws.Cells.Item(138,5).Value = 459
$ws.Cells.Item(138,6).Value = 0
$ws.Cells.Item(138,7).Value = 0
$ws.Cells.Item(138,8).Value = 17

# Row 142: Islandia
$ws.Cells.Item(142,1).Value = "Islandia"
$ws.Cells.Item(142,2).Value = 4230
$ws.Cells.Item(142,3).Value = 37
$ws.Cells.Item(142,4).Value = 3013
$ws.Cells.Item(142,5).Value = 1206
$ws.Cells.Item(142,6).Value = 0
$ws.Cells.Item(142,7).Value = 0
$ws.Cells.Item(142,8).Value = 11

# Row 168: Vietnam
$ws.Cells.Item(168,1).Value = "Vietnam"
$ws.Cells.Item(168,2).Value = 1144
$ws.Cells.Item(168,3).Value = 3
$ws.Cells.Item(168,4).Value = 1046
$ws.Cells.Item(168,5).Value = 63
$ws.Cells.Item(168,6).Value = 0
$ws.Cells.Item(168,7).Value = 0
$ws.Cells.Item(168,8).Value = 35

# Row 175: Gibraltar
$ws.Cells.Item(175,1).Value = "Gibraltar"
$ws.Cells.Item(175,2).Value = 621
$ws.Cells.Item(175,3).Value = 13
$ws.Cells.Item(175,4).Value = 481
$ws.Cells.Item(175,5).Value = 140
$ws.Cells.Item(175,6).Value = 0
$ws.Cells.Item(175,7).Value = 0
$ws.Cells.Item(175,8).Value = 0

# Row 190: Liechtenstein
$ws.Cells.Item(190,1).Value = "Liechtenstein"
$ws.Cells.Item(190,2).Value = 252
$ws.Cells.Item(190,3).Value = 17
$ws.Cells.Item(190,4).Value = 145
$ws.Cells.Item(190,5).Value = 106
$ws.Cells.Item(190,6).Value = 0
$ws.Cells.Item(190,7).Value = 0
$ws.Cells.Item(190,8).Value = 1

# Row 197: Antigua y Barbuda
$ws.Cells.Item(197,1).Value = "Antigua y Barbuda"
$ws.Cells.Item(197,2).Value = 122
$ws.Cells.Item(197,3).Value = 3
$ws.Cells.Item(197,4).Value = 101
$ws.Cells.Item(197,5).Value = 18
$ws.Cells.Item(197,6).Value = 0
$ws.Cells.Item(197,7).Value = 0
$ws.Cells.Item(197,8).Value = 3
